$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "SEO"
$ws.Range("B3").Value = "Pas de meta description"
$ws.Range("C3").Value = "Content vide"
$ws.Range("D3").Value = "Mettre une description"

$ws.Range("A4").Value = "SEO"
$ws.Range("B4").Value = "police trop petite "
$ws.Range("C4").Value = "- de 12px sur 60% du site"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").WrapText = $false
$ws.Range("D4").Value = "augmenter en rem ou en %"

$ws.Range("A5").Value = "SEO"
$ws.Range("B5").Value = "cibles tactiles trop petites"
$ws.Range("C5").Value = "inférieur à 48px par 48 px"
$ws.Range("D5").Value = "augmenter à 48 48px "

$ws.Range("A6").Value = "accessibilité"
$ws.Range("B6").Value = "couleurs d'arrière-plan et de premier plan n'ont pas un rapport de contraste suffisant"
$ws.Range("C6").Value = "pas un rapport de contraste suffisamment élevé"
$ws.Range("D6").Value = "Le texte de 18 points ou 14 points en gras nécessite un rapport de contraste de 3: 1."

$ws.Range("D7").Value = "Tout autre texte a besoin d'un rapport de contraste de 4,5: 1."

$ws.Range("A8").Value = "accessibilité"
$ws.Range("A8").Font.ThemeColor = 1
$ws.Range("A8").WrapText = $false
$ws.Range("B8").Value = "Les éléments d'en-tête ne sont pas dans un ordre séquentiel décroissant"
$ws.Range("C8").Value = "Échec de l'audit des niveaux de titre du phare "
$ws.Range("D8").Value = "utiliser les h1 h2 h3 etc,,,"
$ws.Range("D8").Font.ThemeColor = 1
$ws.Range("D8").WrapText = $false

$ws.Range("A9").Value = "accessibilité"
$ws.Range("A9").Font.ThemeColor = 1
$ws.Range("A9").WrapText = $false
